$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column C in the target range is formatted as Text so that
# numeric-looking labels (e.g. "0", "24", "9", "25") are stored as text,
# matching the source data (t="inlineStr") instead of being auto-converted
# to numbers by Excel.
$ws.Range("C105:C147").NumberFormat = "@"

# Build the new data block for rows 105-147 (row 105 is updated in place,
# rows 106-147 are newly appended).
$arr = New-Object 'object[,]' 43,3
$arr[0,0] = '08/11/2023 10'
$arr[0,1] = 45.85
$arr[0,2] = '0'
$arr[1,0] = '08/11/2023 11'
$arr[1,1] = 44.82
$arr[1,2] = '0'
$arr[2,0] = '08/11/2023 12'
$arr[2,1] = 46.16
$arr[2,2] = '0'
$arr[3,0] = '08/11/2023 13'
$arr[3,1] = 48.49
$arr[3,2] = '0'
$arr[4,0] = '08/11/2023 14'
$arr[4,1] = 51.24
$arr[4,2] = '0'
$arr[5,0] = '08/11/2023 15'
$arr[5,1] = 54.65
$arr[5,2] = '0'
$arr[6,0] = '08/11/2023 16'
$arr[6,1] = 53
$arr[6,2] = '0'
$arr[7,0] = '08/11/2023 17'
$arr[7,1] = 51.2
$arr[7,2] = '0'
$arr[8,0] = '08/11/2023 18'
$arr[8,1] = 49.88
$arr[8,2] = '0'
$arr[9,0] = '08/11/2023 19'
$arr[9,1] = 55.95
$arr[9,2] = '0'
$arr[10,0] = '08/11/2023 20'
$arr[10,1] = 239.28
$arr[10,2] = '0'
$arr[11,0] = '08/11/2023 21'
$arr[11,1] = 592.0700000000001
$arr[11,2] = '0'
$arr[12,0] = '08/11/2023 22'
$arr[12,1] = 582.3
$arr[12,2] = '0'
$arr[13,0] = '08/11/2023 23'
$arr[13,1] = 430.45
$arr[13,2] = '0'
$arr[14,0] = '08/11/2023 24'
$arr[14,1] = 320.9
$arr[14,2] = '0'
$arr[15,0] = '08/12/2023 01'
$arr[15,1] = 491
$arr[15,2] = '0'
$arr[16,0] = '08/12/2023 02'
$arr[16,1] = 341.07
$arr[16,2] = '0'
$arr[17,0] = '08/12/2023 03'
$arr[17,1] = 90.7
$arr[17,2] = '0'
$arr[18,0] = '08/12/2023 04'
$arr[18,1] = -1
$arr[18,2] = '0'
$arr[19,0] = '08/14/2023 09'
$arr[19,1] = 55.03
$arr[19,2] = '0'
$arr[20,0] = '08/14/2023 10'
$arr[20,1] = 64.2
$arr[20,2] = '0'
$arr[21,0] = '08/14/2023 11'
$arr[21,1] = 77.76000000000001
$arr[21,2] = '0'
$arr[22,0] = '08/14/2023 12'
$arr[22,1] = 218.46
$arr[22,2] = '0'
$arr[23,0] = '08/14/2023 19'
$arr[23,1] = 491.44
$arr[23,2] = '24'
$arr[24,0] = '08/14/2023 20'
$arr[24,1] = 710.03
$arr[24,2] = '0'
$arr[25,0] = '08/14/2023 21'
$arr[25,1] = 559.89
$arr[25,2] = '0'
$arr[26,0] = '08/14/2023 22'
$arr[26,1] = 349.73
$arr[26,2] = '0'
$arr[27,0] = '08/14/2023 23'
$arr[27,1] = 233.76
$arr[27,2] = '0'
$arr[28,0] = '08/14/2023 24'
$arr[28,1] = 173.86
$arr[28,2] = '0'
$arr[29,0] = '08/15/2023 01'
$arr[29,1] = 142.91
$arr[29,2] = '0'
$arr[30,0] = '08/15/2023 02'
$arr[30,1] = 68.98999999999999
$arr[30,2] = '0'
$arr[31,0] = '08/15/2023 03'
$arr[31,1] = 53.17
$arr[31,2] = '0'
$arr[32,0] = '08/15/2023 04'
$arr[32,1] = 49.28
$arr[32,2] = '0'
$arr[33,0] = '08/15/2023 05'
$arr[33,1] = 59.62
$arr[33,2] = '0'
$arr[34,0] = '08/15/2023 06'
$arr[34,1] = 155.44
$arr[34,2] = '0'
$arr[35,0] = '08/15/2023 07'
$arr[35,1] = 273.66
$arr[35,2] = '0'
$arr[36,0] = '08/15/2023 08'
$arr[36,1] = 71.54000000000001
$arr[36,2] = '0'
$arr[37,0] = '08/15/2023 09'
$arr[37,1] = 320.14
$arr[37,2] = '0'
$arr[38,0] = '08/15/2023 10'
$arr[38,1] = 111.68
$arr[38,2] = '0'
$arr[39,0] = '08/15/2023 11'
$arr[39,1] = 101.28
$arr[39,2] = '9'
$arr[40,0] = '08/15/2023 12'
$arr[40,1] = 117.42
$arr[40,2] = '25'
$arr[41,0] = '08/15/2023 13'
$arr[41,1] = 269.41
$arr[41,2] = '25'
$arr[42,0] = '08/15/2023 14'
$arr[42,1] = -1
$arr[42,2] = '-'

$ws.Range("A105:C147").Value = $arr

# Update the sheet view to reflect the new scroll position / selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 103
$win.ScrollColumn = 1
$ws.Range("D126").Select()
